# Update cryptocurrency price/volume data per latest GitHub Actions scrape run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $c = $ws.Range("D2")
    $c.NumberFormat = "@"
    $c.Value = '42.991.34'
    $c.Style = "Normal"
    $ws.Range("E2").Value = '  -0.31%  '
    $c = $ws.Range("D3")
    $c.NumberFormat = "@"
    $c.Value = '2.553.70'
    $c.Style = "Normal"
    $ws.Range("E3").Value = '  -0.02%  '
    $ws.Range("E4").Value = '  -0.14%  '
    $c = $ws.Range("D5")
    $c.NumberFormat = "@"
    $c.Value = '304.96'
    $c.Style = "Normal"
    $ws.Range("E5").Value = '  +1.73%  '
    $c = $ws.Range("D6")
    $c.NumberFormat = "@"
    $c.Value = '98.30'
    $c.Style = "Normal"
    $ws.Range("E6").Value = '  +5.85%  '
    $c = $ws.Range("D7")
    $c.NumberFormat = "@"
    $c.Value = '0.577'
    $c.Style = "Normal"
    $ws.Range("E7").Value = '  +0.30%  '
    $ws.Range("E8").Value = '  +0.09%  '
    $c = $ws.Range("D9")
    $c.NumberFormat = "@"
    $c.Value = '0.547'
    $c.Style = "Normal"
    $ws.Range("E9").Value = '  -0.62%  '
    $c = $ws.Range("D10")
    $c.NumberFormat = "@"
    $c.Value = '37.10'
    $c.Style = "Normal"
    $ws.Range("E10").Value = '  +3.31%  '
    $c = $ws.Range("D11")
    $c.NumberFormat = "@"
    $c.Value = '0.0822'
    $c.Style = "Normal"
    $ws.Range("E11").Value = '  +2.00%  '
    $ws.Range("B12").Value = 'TRON'
    $ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    $c = $ws.Range("D12")
    $c.NumberFormat = "@"
    $c.Value = '0.116'
    $c.Style = "Normal"
    $ws.Range("E12").Value = '  +5.79%  '
    $ws.Range("B13").Value = 'Polkadot'
    $ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    $c = $ws.Range("D13")
    $c.NumberFormat = "@"
    $c.Value = '7.64'
    $c.Style = "Normal"
    $ws.Range("E13").Value = '  -1.44%  '
    $c = $ws.Range("D14")
    $c.NumberFormat = "@"
    $c.Value = '2.945.44'
    $c.Style = "Normal"
    $ws.Range("E14").Value = '  +0.01%  '
    $c = $ws.Range("D15")
    $c.NumberFormat = "@"
    $c.Value = '2.620.93'
    $c.Style = "Normal"
    $ws.Range("E15").Value = '  +1.58%  '
    $c = $ws.Range("D16")
    $c.NumberFormat = "@"
    $c.Value = '15.00'
    $c.Style = "Normal"
    $ws.Range("E16").Value = '  +6.22%  '
    $c = $ws.Range("D17")
    $c.NumberFormat = "@"
    $c.Value = '0.879'
    $c.Style = "Normal"
    $ws.Range("E17").Value = '  +0.62%  '
    $c = $ws.Range("D18")
    $c.NumberFormat = "@"
    $c.Value = '43.014.72'
    $c.Style = "Normal"
    $ws.Range("E18").Value = '  -0.32%  '
    $c = $ws.Range("D19")
    $c.NumberFormat = "@"
    $c.Value = '13.76'
    $c.Style = "Normal"
    $ws.Range("E19").Value = '  +4.47%  '
    $c = $ws.Range("D20")
    $c.NumberFormat = "@"
    $c.Value = '0.0₃0996'
    $c.Style = "Normal"
    $ws.Range("E20").Value = '  +1.69%  '
    $c = $ws.Range("D21")
    $c.NumberFormat = "@"
    $c.Value = '6.61'
    $c.Style = "Normal"
    $ws.Range("E21").Value = '  -0.26%  '
    $c = $ws.Range("D22")
    $c.NumberFormat = "@"
    $c.Value = '72.00'
    $c.Style = "Normal"
    $ws.Range("E22").Value = '  +0.06%  '
    $c = $ws.Range("D23")
    $c.NumberFormat = "@"
    $c.Value = '254.42'
    $c.Style = "Normal"
    $ws.Range("E23").Value = '  -1.90%  '
    $ws.Range("E24").Value = '  +0.87%  '
    $c = $ws.Range("D25")
    $c.NumberFormat = "@"
    $c.Value = '2.09'
    $c.Style = "Normal"
    $ws.Range("E25").Value = '  -2.27%  '
    $c = $ws.Range("D26")
    $c.NumberFormat = "@"
    $c.Value = '28.12'
    $c.Style = "Normal"
    $ws.Range("E26").Value = '  -4.89%  '
    $ws.Range("E27").Value = '  -0.07%  '
    $c = $ws.Range("D28")
    $c.NumberFormat = "@"
    $c.Value = '10.20'
    $c.Style = "Normal"
    $ws.Range("E28").Value = '  +1.25%  '
    $c = $ws.Range("D29")
    $c.NumberFormat = "@"
    $c.Value = '38.12'
    $c.Style = "Normal"
    $ws.Range("E29").Value = '  +0.73%  '
    $ws.Range("E30").Value = '  -1.07%  '
    $c = $ws.Range("D31")
    $c.NumberFormat = "@"
    $c.Value = '6.16'
    $c.Style = "Normal"
    $ws.Range("E31").Value = '  +2.82%  '
    $c = $ws.Range("D32")
    $c.NumberFormat = "@"
    $c.Value = '158.79'
    $c.Style = "Normal"
    $ws.Range("E32").Value = '  +2.77%  '
    $c = $ws.Range("D33")
    $c.NumberFormat = "@"
    $c.Value = '2.16'
    $c.Style = "Normal"
    $ws.Range("E33").Value = '  -0.90%  '
    $ws.Range("E34").Value = '  -0.46%  '
    $c = $ws.Range("D35")
    $c.NumberFormat = "@"
    $c.Value = '19.46'
    $c.Style = "Normal"
    $ws.Range("E35").Value = '  +14.79%  '
    $c = $ws.Range("D36")
    $c.NumberFormat = "@"
    $c.Value = '0.0804'
    $c.Style = "Normal"
    $ws.Range("E36").Value = '  +1.09%  '
    $ws.Range("E37").Value = '  -2.03%  '
    $ws.Range("E38").Value = '  +0.65%  '
    $c = $ws.Range("D39")
    $c.NumberFormat = "@"
    $c.Value = '25.84'
    $c.Style = "Normal"
    $ws.Range("E39").Value = '  +10.04%  '
    $ws.Range("E40").Value = '  -0.12%  '
    $ws.Range("E41").Value = '  +31.89%  '
    $c = $ws.Range("D42")
    $c.NumberFormat = "@"
    $c.Value = '3.44'
    $c.Style = "Normal"
    $ws.Range("E42").Value = '  -0.57%  '
    $ws.Range("E43").Value = '  +0.09%  '
    $ws.Range("E44").Value = '  -1.68%  '
    $c = $ws.Range("D45")
    $c.NumberFormat = "@"
    $c.Value = '2.082.19'
    $c.Style = "Normal"
    $ws.Range("E45").Value = '  +0.16%  '
    $ws.Range("E46").Value = '  +0.03%  '
    $ws.Range("E47").Value = '  +1.14%  '
    $c = $ws.Range("D48")
    $c.NumberFormat = "@"
    $c.Value = '8.99'
    $c.Style = "Normal"
    $ws.Range("E48").Value = '  +0.10%  '
    $c = $ws.Range("D49")
    $c.NumberFormat = "@"
    $c.Value = '2.802.75'
    $c.Style = "Normal"
    $ws.Range("E49").Value = '  +0.01%  '
    $c = $ws.Range("D50")
    $c.NumberFormat = "@"
    $c.Value = '74.94'
    $c.Style = "Normal"
    $ws.Range("E50").Value = '  +8.27%  '
    $c = $ws.Range("D51")
    $c.NumberFormat = "@"
    $c.Value = '103.39'
    $c.Style = "Normal"
    $ws.Range("E51").Value = '  -1.73%  '
